$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "x" column (D) for the relevant battle rows (Erledigt list)
$ws.Range("D3").Value = "x"
$ws.Range("D4").Value = "x"
$ws.Range("D6").Value = "x"
$ws.Range("D25").Value = "x"

# Move the active selection to A7
$ws.Range("A7").Select()
